$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column D width ---
$ws.Columns.Item(4).ColumnWidth = 29.1

# --- 2. Insert 2 new rows after row 17 (row17 content/style stays put at row 17 for now) ---
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# --- 3. Stash row17's "last-row" (bottom border) formatting to the far right of row 17 ---
$ws.Range("B17:J17").Copy()
$ws.Range("AB17:AJ17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Apply "normal row" formatting (from row16) onto rows 17 and 18 ---
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Apply stashed "last-row" formatting onto row 19 ---
$ws.Range("AB17:AJ17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 6. Clear scratch area ---
$ws.Range("AB17:AJ17").Clear()

# --- 7. Fill in values ---
# Row 16: CC | 1053122836 | YOLADIS PAJARO MATA | 2507 | 56940 | 877803
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1053122836"
$ws.Range("D16").Value = "YOLADIS PAJARO MATA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 877803

# Row 17: CC | 1053122836 | YOLADIS PAJARO MATA | 2506 | 56940 | 877803
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1053122836"
$ws.Range("D17").Value = "YOLADIS PAJARO MATA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 877803

# Row 18: CC | 1053122836 | YOLADIS PAJARO MATA | 2505 | 56940 | 877803
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1053122836"
$ws.Range("D18").Value = "YOLADIS PAJARO MATA"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 877803

# Row 19: CC | 1143406582 | ZULENA MARIA ARRIETA PUELLO | 2507 | 96000 | 2400000
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143406582"
$ws.Range("D19").Value = "ZULENA MARIA ARRIETA PUELLO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 96000
$ws.Range("G19").Value = 2400000

# --- 8. Update summary values ---
$ws.Range("E11").Value = 266820
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3

Write-Host "all done"
